$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 1.04
$ws.Range("K2").Value = 13
$ws.Range("L2").Value = 1.22
$ws.Range("M2").Value = 4.33
$ws.Range("N2").Value = 1.73
$ws.Range("O2").Value = 2.1
$ws.Range("W2").Value = 21
$ws.Range("AA2").Value = 7
$ws.Range("AB2").Value = 13
$ws.Range("AF2").Value = 21
$ws.Range("AI2").Value = 29

# Row 4
$ws.Range("G4").Value = 1.25
$ws.Range("I4").Value = 9
$ws.Range("R4").Value = 1.91
$ws.Range("S4").Value = 1.91
$ws.Range("AE4").Value = 26
$ws.Range("AF4").Value = 51
$ws.Range("AG4").Value = 26
$ws.Range("AH4").Value = 126

# Row 5
$ws.Range("K5").Value = 7.5
$ws.Range("R5").Value = 1.91
$ws.Range("S5").Value = 1.8
$ws.Range("T5").Value = 6.5
$ws.Range("U5").Value = 9.5
$ws.Range("Y5").Value = 34
$ws.Range("Z5").Value = 7.5
$ws.Range("AD5").Value = 351
$ws.Range("AE5").Value = 9.5
$ws.Range("AI5").Value = 34

# Row 6
$ws.Range("N6").Value = 1.92
$ws.Range("O6").Value = 1.82

# Row 7
$ws.Range("N7").Value = 1.82
$ws.Range("O7").Value = 1.92

# Row 8
$ws.Range("N8").Value = 1.79
$ws.Range("O8").Value = 1.94

# Row 13
$ws.Range("G13").Value = 2.4
$ws.Range("H13").Value = 3
$ws.Range("I13").Value = 2.95
$ws.Range("P13").Value = 1.53
$ws.Range("Q13").Value = 2.18
$ws.Range("R13").Value = 1.98
$ws.Range("T13").Value = 6.4
$ws.Range("U13").Value = 10.5
$ws.Range("V13").Value = 9.75
$ws.Range("W13").Value = 25
$ws.Range("Y13").Value = 40
$ws.Range("Z13").Value = 6.9
$ws.Range("AA13").Value = 5.9
$ws.Range("AB13").Value = 17.5
$ws.Range("AH13").Value = 37
$ws.Range("AI13").Value = 32

# Row 16
$ws.Range("G16").Value = 2.25
$ws.Range("H16").Value = 2.95
$ws.Range("I16").Value = 3.35
$ws.Range("J16").Value = 1.08
$ws.Range("K16").Value = 6.4
$ws.Range("L16").Value = 1.36
$ws.Range("M16").Value = 2.9
$ws.Range("N16").Value = 2.05
$ws.Range("O16").Value = 1.7
$ws.Range("P16").Value = 1.4
$ws.Range("Q16").Value = 2.72
$ws.Range("R16").Value = 1.78
$ws.Range("S16").Value = 1.93
$ws.Range("T16").Value = 7.2
$ws.Range("U16").Value = 10.75
$ws.Range("W16").Value = 23
$ws.Range("X16").Value = 18.5
$ws.Range("Y16").Value = 29
$ws.Range("Z16").Value = 6.4
$ws.Range("AA16").Value = 5.7
$ws.Range("AC16").Value = 65
$ws.Range("AD16").Value = 500
$ws.Range("AE16").Value = 9.25
$ws.Range("AF16").Value = 18
$ws.Range("AG16").Value = 11.25
$ws.Range("AI16").Value = 32
$ws.Range("AJ16").Value = 37

# Row 17
$ws.Range("G17").Value = 3.25
$ws.Range("H17").Value = 2.77
$ws.Range("I17").Value = 2.42
$ws.Range("K17").Value = 5.5
$ws.Range("R17").Value = 1.93
$ws.Range("S17").Value = 1.78
$ws.Range("T17").Value = 7.9
$ws.Range("U17").Value = 16
$ws.Range("V17").Value = 11.25
$ws.Range("W17").Value = 45
$ws.Range("X17").Value = 32
$ws.Range("Z17").Value = 5.5
$ws.Range("AE17").Value = 6.3
$ws.Range("AF17").Value = 10.75
$ws.Range("AH17").Value = 26
$ws.Range("AI17").Value = 24

# Row 18
$ws.Range("G18").Value = 1.65
$ws.Range("H18").Value = 4.05
$ws.Range("I18").Value = 4.15
$ws.Range("N18").Value = 1.55
$ws.Range("O18").Value = 2.15
$ws.Range("R18").Value = 1.57
$ws.Range("S18").Value = 2.1
$ws.Range("T18").Value = 9.25
$ws.Range("X18").Value = 12
$ws.Range("Y18").Value = 21
$ws.Range("Z18").Value = 15
$ws.Range("AA18").Value = 8.25
$ws.Range("AB18").Value = 14
$ws.Range("AC18").Value = 50
$ws.Range("AD18").Value = 300
$ws.Range("AE18").Value = 15.5
$ws.Range("AF18").Value = 26
$ws.Range("AG18").Value = 14
$ws.Range("AH18").Value = 65

# Row 21
$ws.Range("G21").Value = 2.45
$ws.Range("AE21").Value = 9.5

# Row 22
$ws.Range("I22").Value = 4.33

# Row 23
$ws.Range("G23").Value = 2.05
$ws.Range("H23").Value = 3.1
$ws.Range("I23").Value = 4

# Row 24
$ws.Range("G24").Value = 1.95
$ws.Range("H24").Value = 3.1
$ws.Range("I24").Value = 4.75
$ws.Range("N24").Value = 2.08
$ws.Range("O24").Value = 1.73
$ws.Range("AF24").Value = 21

# Row 30
$ws.Range("G30").Value = 1.75
$ws.Range("I30").Value = 4.1
$ws.Range("L30").Value = 1.3
$ws.Range("P30").Value = 1.38
$ws.Range("Q30").Value = 2.45
$ws.Range("R30").Value = 1.84
$ws.Range("S30").Value = 1.86
$ws.Range("T30").Value = 5.7
$ws.Range("U30").Value = 6.8
$ws.Range("V30").Value = 7
$ws.Range("W30").Value = 11.5
$ws.Range("Z30").Value = 9.25
$ws.Range("AC30").Value = 50
$ws.Range("AE30").Value = 9.5
$ws.Range("AF30").Value = 18.5
$ws.Range("AG30").Value = 11.25
$ws.Range("AH30").Value = 50
$ws.Range("AI30").Value = 30

# Row 35
$ws.Range("L35").Value = 1.25
$ws.Range("M35").Value = 3.75
$ws.Range("N35").Value = 1.8
$ws.Range("O35").Value = 2
